$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26 (shifts existing rows 26-74 down to 27-75)
$ws.Rows.Item(26).Insert()

# Populate the new row 26 with the new weekly record
$ws.Range("A26").Value = 1
$ws.Range("B26").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C26").Value = "Arica y Parinacota"
$ws.Range("D26").Value = 44791
$ws.Range("E26").Value = 15
$ws.Range("F26").Value = 100112012
$ws.Range("G26").Value = "Espinaca"
$ws.Range("H26").Value = "Sin especificar"
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 300
$ws.Range("K26").Value = 2000
$ws.Range("L26").Value = 2500
$ws.Range("M26").Value = 2250
$ws.Range("N26").Value = "`$/atado 2,5 a 3 kilos"
$ws.Range("O26").Value = "Región de Arica y Parinacota"
$ws.Range("P26").Value = 750
$ws.Range("Q26").Value = 3
$ws.Range("R26").Value = "Hortaliza"
